$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 2).Value = 6726054
$ws.Cells.Item(6, 2).Value = 6221786
$ws.Cells.Item(5, 5).Value = "FK Aktobe"
$ws.Cells.Item(6, 5).Value = "Ordabasy"
$ws.Cells.Item(5, 6).Value = "Shakhter Karagandy"
$ws.Cells.Item(6, 6).Value = "FK Atyrau"
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(5, 10).Value = 1.727
$ws.Cells.Item(6, 10).Value = 1.571
$ws.Cells.Item(5, 11).Value = 4
$ws.Cells.Item(6, 11).Value = 3.8
$ws.Cells.Item(5, 12).Value = 3.5
$ws.Cells.Item(6, 12).Value = 4.75
$ws.Cells.Item(5, 13).Value = 1.333
$ws.Cells.Item(6, 13).Value = 1.5
$ws.Cells.Item(5, 14).Value = 4.75
$ws.Cells.Item(6, 14).Value = 4
$ws.Cells.Item(5, 15).Value = 6.5
$ws.Cells.Item(6, 15).Value = 5.25
$ws.Cells.Item(5, 16).Value = -1.5
$ws.Cells.Item(6, 16).Value = -1
$ws.Cells.Item(5, 17).Value = 1.975
$ws.Cells.Item(6, 17).Value = 1.875
$ws.Cells.Item(5, 18).Value = 1.825
$ws.Cells.Item(6, 18).Value = 1.925
$ws.Cells.Item(5, 19).Value = 3
$ws.Cells.Item(6, 19).Value = 2.5
$ws.Cells.Item(5, 20).Value = 1.975
$ws.Cells.Item(6, 20).Value = 1.9
$ws.Cells.Item(5, 21).Value = 1.825
$ws.Cells.Item(6, 21).Value = 1.9
$ws.Cells.Item(5, 22).Value = 0.333
$ws.Cells.Item(6, 22).Value = 0.5
$ws.Cells.Item(5, 25).Value = 0.9750000000000001
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(5, 26).Value = -1
$ws.Cells.Item(6, 26).Value = 0
$ws.Cells.Item(5, 27).Value = -1
$ws.Cells.Item(6, 27).Value = 0.8999999999999999
$ws.Cells.Item(5, 28).Value = 0.825
$ws.Cells.Item(6, 28).Value = -1
$ws.Cells.Item(16, 2).Value = 6221698
$ws.Cells.Item(17, 2).Value = 6221693
$ws.Cells.Item(16, 6).Value = "FK Aktobe"
$ws.Cells.Item(17, 6).Value = "Shakhter Karagandy"
$ws.Cells.Item(16, 8).Value = 2
$ws.Cells.Item(17, 8).Value = 3
$ws.Cells.Item(16, 10).Value = 4.333
$ws.Cells.Item(17, 10).Value = 2
$ws.Cells.Item(16, 11).Value = 3.5
$ws.Cells.Item(17, 11).Value = 3.4
$ws.Cells.Item(16, 12).Value = 1.666
$ws.Cells.Item(17, 12).Value = 3.1
$ws.Cells.Item(16, 13).Value = 4.2
$ws.Cells.Item(17, 13).Value = 2.2
$ws.Cells.Item(16, 14).Value = 3.4
$ws.Cells.Item(17, 14).Value = 3.3
$ws.Cells.Item(16, 15).Value = 1.7
$ws.Cells.Item(17, 15).Value = 2.8
$ws.Cells.Item(16, 16).Value = 0.75
$ws.Cells.Item(17, 16).Value = -0.25
$ws.Cells.Item(16, 17).Value = 1.825
$ws.Cells.Item(17, 17).Value = 1.95
$ws.Cells.Item(16, 18).Value = 1.975
$ws.Cells.Item(17, 18).Value = 1.85
$ws.Cells.Item(16, 20).Value = 1.925
$ws.Cells.Item(17, 20).Value = 1.85
$ws.Cells.Item(16, 21).Value = 1.875
$ws.Cells.Item(17, 21).Value = 1.95
$ws.Cells.Item(16, 24).Value = 0.7
$ws.Cells.Item(17, 24).Value = 1.8
$ws.Cells.Item(16, 25).Value = -0.5
$ws.Cells.Item(17, 25).Value = -1
$ws.Cells.Item(16, 26).Value = 0.4875
$ws.Cells.Item(17, 26).Value = 0.8500000000000001
$ws.Cells.Item(16, 27).Value = 0.925
$ws.Cells.Item(17, 27).Value = 0.8500000000000001
$ws.Cells.Item(25, 2).Value = 6221699
$ws.Cells.Item(26, 2).Value = 6221703
$ws.Cells.Item(25, 6).Value = "Kaisar Kyzylorda"
$ws.Cells.Item(26, 6).Value = "FK Aktobe"
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 2
$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(25, 9).Value = "D"
$ws.Cells.Item(26, 9).Value = "A"
$ws.Cells.Item(25, 10).Value = 3.1
$ws.Cells.Item(26, 10).Value = 3.6
$ws.Cells.Item(25, 11).Value = 3.2
$ws.Cells.Item(26, 11).Value = 3.5
$ws.Cells.Item(25, 12).Value = 2.1
$ws.Cells.Item(26, 12).Value = 1.8
$ws.Cells.Item(25, 13).Value = 2.1
$ws.Cells.Item(26, 13).Value = 3.1
$ws.Cells.Item(25, 14).Value = 3.1
$ws.Cells.Item(26, 14).Value = 3.5
$ws.Cells.Item(25, 15).Value = 3.2
$ws.Cells.Item(26, 15).Value = 1.909
$ws.Cells.Item(25, 16).Value = -0.25
$ws.Cells.Item(26, 16).Value = 0.5
$ws.Cells.Item(25, 17).Value = 1.85
$ws.Cells.Item(26, 17).Value = 1.825
$ws.Cells.Item(25, 18).Value = 1.95
$ws.Cells.Item(26, 18).Value = 1.975
$ws.Cells.Item(25, 19).Value = 2.25
$ws.Cells.Item(26, 19).Value = 2.5
$ws.Cells.Item(25, 20).Value = 1.975
$ws.Cells.Item(26, 20).Value = 1.75
$ws.Cells.Item(25, 21).Value = 1.725
$ws.Cells.Item(26, 21).Value = 1.95
$ws.Cells.Item(25, 23).Value = 2.1
$ws.Cells.Item(26, 23).Value = -1
$ws.Cells.Item(25, 24).Value = -1
$ws.Cells.Item(26, 24).Value = 0.909
$ws.Cells.Item(25, 25).Value = -0.5
$ws.Cells.Item(26, 25).Value = -1
$ws.Cells.Item(25, 26).Value = 0.475
$ws.Cells.Item(26, 26).Value = 0.9750000000000001
$ws.Cells.Item(25, 27).Value = 0.9750000000000001
$ws.Cells.Item(26, 27).Value = -1
$ws.Cells.Item(25, 28).Value = -1
$ws.Cells.Item(26, 28).Value = 0.95
$ws.Cells.Item(50, 2).Value = 7055064
$ws.Cells.Item(51, 2).Value = 6221723
$ws.Cells.Item(50, 6).Value = "Ordabasy"
$ws.Cells.Item(51, 6).Value = "FK Aktobe"
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 3
$ws.Cells.Item(51, 8).Value = 1
$ws.Cells.Item(50, 9).Value = "A"
$ws.Cells.Item(51, 9).Value = "D"
$ws.Cells.Item(50, 10).Value = 6.5
$ws.Cells.Item(51, 10).Value = 3.75
$ws.Cells.Item(50, 11).Value = 4.5
$ws.Cells.Item(51, 11).Value = 3.3
$ws.Cells.Item(50, 12).Value = 1.363
$ws.Cells.Item(51, 12).Value = 1.833
$ws.Cells.Item(50, 13).Value = 4.2
$ws.Cells.Item(51, 13).Value = 2.6
$ws.Cells.Item(50, 14).Value = 4
$ws.Cells.Item(51, 14).Value = 3.1
$ws.Cells.Item(50, 15).Value = 1.6
$ws.Cells.Item(51, 15).Value = 2.5
$ws.Cells.Item(50, 16).Value = 1
$ws.Cells.Item(51, 16).Value = 0
$ws.Cells.Item(50, 17).Value = 1.725
$ws.Cells.Item(51, 17).Value = 1.925
$ws.Cells.Item(50, 18).Value = 2.075
$ws.Cells.Item(51, 18).Value = 1.875
$ws.Cells.Item(50, 19).Value = 2.25
$ws.Cells.Item(51, 19).Value = 2.5
$ws.Cells.Item(50, 20).Value = 1.875
$ws.Cells.Item(51, 20).Value = 1.9
$ws.Cells.Item(50, 21).Value = 1.925
$ws.Cells.Item(51, 21).Value = 1.9
$ws.Cells.Item(50, 23).Value = -1
$ws.Cells.Item(51, 23).Value = 2.1
$ws.Cells.Item(50, 24).Value = 0.6000000000000001
$ws.Cells.Item(51, 24).Value = -1
$ws.Cells.Item(50, 27).Value = 0.875
$ws.Cells.Item(51, 27).Value = -1
$ws.Cells.Item(50, 28).Value = -1
$ws.Cells.Item(51, 28).Value = 0.8999999999999999
$ws.Cells.Item(63, 2).Value = 6221732
$ws.Cells.Item(64, 2).Value = 6221729
$ws.Cells.Item(63, 6).Value = "FK Aktobe"
$ws.Cells.Item(64, 6).Value = "FC Astana"
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(63, 9).Value = "D"
$ws.Cells.Item(64, 9).Value = "H"
$ws.Cells.Item(63, 10).Value = 3.3
$ws.Cells.Item(64, 10).Value = 2.9
$ws.Cells.Item(63, 11).Value = 3.25
$ws.Cells.Item(64, 11).Value = 3.3
$ws.Cells.Item(63, 12).Value = 2
$ws.Cells.Item(64, 12).Value = 2.15
$ws.Cells.Item(63, 13).Value = 3.5
$ws.Cells.Item(64, 13).Value = 1.75
$ws.Cells.Item(63, 14).Value = 3.25
$ws.Cells.Item(64, 14).Value = 3.5
$ws.Cells.Item(63, 15).Value = 1.909
$ws.Cells.Item(64, 15).Value = 4
$ws.Cells.Item(63, 16).Value = 0.5
$ws.Cells.Item(64, 16).Value = -0.5
$ws.Cells.Item(63, 17).Value = 1.775
$ws.Cells.Item(64, 17).Value = 1.8
$ws.Cells.Item(63, 18).Value = 2.025
$ws.Cells.Item(64, 18).Value = 2
$ws.Cells.Item(63, 19).Value = 2.25
$ws.Cells.Item(64, 19).Value = 2.5
$ws.Cells.Item(63, 22).Value = -1
$ws.Cells.Item(64, 22).Value = 0.75
$ws.Cells.Item(63, 23).Value = 2.25
$ws.Cells.Item(64, 23).Value = -1
$ws.Cells.Item(63, 25).Value = 0.7749999999999999
$ws.Cells.Item(64, 25).Value = 0.8
$ws.Cells.Item(85, 2).Value = 6221809
$ws.Cells.Item(86, 2).Value = 6221743
$ws.Cells.Item(85, 6).Value = "Okzhetpes Kokshetau"
$ws.Cells.Item(86, 6).Value = "Kaisar Kyzylorda"
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 2
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(85, 9).Value = "A"
$ws.Cells.Item(86, 9).Value = "H"
$ws.Cells.Item(85, 10).Value = 2.05
$ws.Cells.Item(86, 10).Value = 3
$ws.Cells.Item(85, 12).Value = 3.3
$ws.Cells.Item(86, 12).Value = 2.2
$ws.Cells.Item(85, 13).Value = 1.55
$ws.Cells.Item(86, 13).Value = 2.9
$ws.Cells.Item(85, 14).Value = 3.5
$ws.Cells.Item(86, 14).Value = 3.1
$ws.Cells.Item(85, 15).Value = 5.25
$ws.Cells.Item(86, 15).Value = 2.25
$ws.Cells.Item(85, 16).Value = -1
$ws.Cells.Item(86, 16).Value = 0.25
$ws.Cells.Item(85, 17).Value = 2
$ws.Cells.Item(86, 17).Value = 1.8
$ws.Cells.Item(85, 18).Value = 1.8
$ws.Cells.Item(86, 18).Value = 2
$ws.Cells.Item(85, 20).Value = 1.875
$ws.Cells.Item(86, 20).Value = 1.95
$ws.Cells.Item(85, 21).Value = 1.925
$ws.Cells.Item(86, 21).Value = 1.85
$ws.Cells.Item(85, 22).Value = -1
$ws.Cells.Item(86, 22).Value = 1.9
$ws.Cells.Item(85, 24).Value = 4.25
$ws.Cells.Item(86, 24).Value = -1
$ws.Cells.Item(85, 25).Value = -1
$ws.Cells.Item(86, 25).Value = 0.8
$ws.Cells.Item(85, 26).Value = 0.8
$ws.Cells.Item(86, 26).Value = -1
$ws.Cells.Item(85, 27).Value = 0.875
$ws.Cells.Item(86, 27).Value = -0.5
$ws.Cells.Item(85, 28).Value = -1
$ws.Cells.Item(86, 28).Value = 0.425
$ws.Cells.Item(88, 2).Value = 6221811
$ws.Cells.Item(89, 2).Value = 6221810
$ws.Cells.Item(88, 6).Value = "Kairat Almaty"
$ws.Cells.Item(89, 6).Value = "FK Maktaaral"
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 1
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(88, 9).Value = "D"
$ws.Cells.Item(89, 9).Value = "H"
$ws.Cells.Item(88, 10).Value = 3.75
$ws.Cells.Item(89, 10).Value = 1.909
$ws.Cells.Item(88, 11).Value = 3.4
$ws.Cells.Item(89, 11).Value = 3.1
$ws.Cells.Item(88, 12).Value = 1.8
$ws.Cells.Item(89, 12).Value = 3.8
$ws.Cells.Item(88, 13).Value = 3.8
$ws.Cells.Item(89, 13).Value = 2
$ws.Cells.Item(88, 14).Value = 3.6
$ws.Cells.Item(89, 14).Value = 3.1
$ws.Cells.Item(88, 15).Value = 1.75
$ws.Cells.Item(89, 15).Value = 3.5
$ws.Cells.Item(88, 16).Value = 0.75
$ws.Cells.Item(89, 16).Value = -0.5
$ws.Cells.Item(88, 17).Value = 1.8
$ws.Cells.Item(89, 17).Value = 2
$ws.Cells.Item(88, 18).Value = 2
$ws.Cells.Item(89, 18).Value = 1.8
$ws.Cells.Item(88, 19).Value = 2.75
$ws.Cells.Item(89, 19).Value = 2
$ws.Cells.Item(88, 20).Value = 1.825
$ws.Cells.Item(89, 20).Value = 1.725
$ws.Cells.Item(88, 22).Value = -1
$ws.Cells.Item(89, 22).Value = 1
$ws.Cells.Item(88, 23).Value = 2.6
$ws.Cells.Item(89, 23).Value = -1
$ws.Cells.Item(88, 25).Value = 0.8
$ws.Cells.Item(89, 25).Value = 1
$ws.Cells.Item(88, 27).Value = -1
$ws.Cells.Item(89, 27).Value = 0
$ws.Cells.Item(88, 28).Value = 0.9750000000000001
$ws.Cells.Item(89, 28).Value = 0
$ws.Cells.Item(98, 2).Value = 6221753
$ws.Cells.Item(99, 2).Value = 6221752
$ws.Cells.Item(98, 6).Value = "Tobol Kostanay"
$ws.Cells.Item(99, 6).Value = "Kaisar Kyzylorda"
$ws.Cells.Item(98, 8).Value = 3
$ws.Cells.Item(99, 8).Value = 1
$ws.Cells.Item(98, 10).Value = 2.75
$ws.Cells.Item(99, 10).Value = 1.833
$ws.Cells.Item(98, 11).Value = 3.1
$ws.Cells.Item(99, 11).Value = 3.2
$ws.Cells.Item(98, 12).Value = 2.375
$ws.Cells.Item(99, 12).Value = 4
$ws.Cells.Item(98, 13).Value = 2.625
$ws.Cells.Item(99, 13).Value = 1.85
$ws.Cells.Item(98, 15).Value = 2.45
$ws.Cells.Item(99, 15).Value = 4
$ws.Cells.Item(98, 16).Value = 0
$ws.Cells.Item(99, 16).Value = -0.5
$ws.Cells.Item(98, 17).Value = 2
$ws.Cells.Item(99, 17).Value = 1.9
$ws.Cells.Item(98, 18).Value = 1.8
$ws.Cells.Item(99, 18).Value = 1.9
$ws.Cells.Item(98, 19).Value = 2.5
$ws.Cells.Item(99, 19).Value = 2
$ws.Cells.Item(98, 20).Value = 1.9
$ws.Cells.Item(99, 20).Value = 1.775
$ws.Cells.Item(98, 21).Value = 1.9
$ws.Cells.Item(99, 21).Value = 2.025
$ws.Cells.Item(98, 24).Value = 1.45
$ws.Cells.Item(99, 24).Value = 3
$ws.Cells.Item(98, 26).Value = 0.8
$ws.Cells.Item(99, 26).Value = 0.8999999999999999
$ws.Cells.Item(98, 27).Value = 0.8999999999999999
$ws.Cells.Item(99, 27).Value = -1
$ws.Cells.Item(98, 28).Value = -1
$ws.Cells.Item(99, 28).Value = 1.025
$ws.Cells.Item(101, 2).Value = 6221755
$ws.Cells.Item(102, 2).Value = 6221754
$ws.Cells.Item(102, 5).Value = "Ordabasy"
$ws.Cells.Item(101, 6).Value = "Zhetysu"
$ws.Cells.Item(102, 6).Value = "FC Astana"
$ws.Cells.Item(101, 7).Value = 2
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 2
$ws.Cells.Item(102, 8).Value = 1
$ws.Cells.Item(101, 10).Value = 1.727
$ws.Cells.Item(102, 10).Value = 3.6
$ws.Cells.Item(101, 12).Value = 4
$ws.Cells.Item(102, 12).Value = 1.8
$ws.Cells.Item(101, 13).Value = 1.444
$ws.Cells.Item(102, 13).Value = 5
$ws.Cells.Item(101, 14).Value = 4
$ws.Cells.Item(102, 14).Value = 1.4
$ws.Cells.Item(101, 15).Value = 6
$ws.Cells.Item(102, 15).Value = 5
$ws.Cells.Item(101, 16).Value = -1.25
$ws.Cells.Item(102, 16).Value = 0.25
$ws.Cells.Item(101, 17).Value = 1.975
$ws.Cells.Item(102, 17).Value = 1.7
$ws.Cells.Item(101, 18).Value = 1.825
$ws.Cells.Item(102, 18).Value = 2.1
$ws.Cells.Item(101, 19).Value = 2.75
$ws.Cells.Item(102, 19).Value = 2.25
$ws.Cells.Item(101, 20).Value = 1.8
$ws.Cells.Item(102, 20).Value = 1.9
$ws.Cells.Item(101, 21).Value = 2
$ws.Cells.Item(102, 21).Value = 1.9
$ws.Cells.Item(101, 23).Value = 3
$ws.Cells.Item(102, 23).Value = 0.3999999999999999
$ws.Cells.Item(101, 25).Value = -1
$ws.Cells.Item(102, 25).Value = 0.35
$ws.Cells.Item(101, 26).Value = 0.825
$ws.Cells.Item(102, 26).Value = -0.5
$ws.Cells.Item(101, 27).Value = 0.8
$ws.Cells.Item(102, 27).Value = -0.5
$ws.Cells.Item(101, 28).Value = -1
$ws.Cells.Item(102, 28).Value = 0.45
$ws.Cells.Item(103, 2).Value = 6221816
$ws.Cells.Item(104, 2).Value = 6221814
$ws.Cells.Item(103, 6).Value = "Ordabasy"
$ws.Cells.Item(104, 6).Value = "FK Maktaaral"
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(104, 8).Value = 1
$ws.Cells.Item(103, 9).Value = "H"
$ws.Cells.Item(104, 9).Value = "D"
$ws.Cells.Item(103, 10).Value = 3.4
$ws.Cells.Item(104, 10).Value = 2.3
$ws.Cells.Item(103, 11).Value = 3.4
$ws.Cells.Item(104, 11).Value = 3.1
$ws.Cells.Item(103, 12).Value = 1.909
$ws.Cells.Item(104, 12).Value = 2.8
$ws.Cells.Item(103, 13).Value = 4.2
$ws.Cells.Item(104, 13).Value = 2.3
$ws.Cells.Item(103, 14).Value = 4
$ws.Cells.Item(104, 14).Value = 3.1
$ws.Cells.Item(103, 15).Value = 1.571
$ws.Cells.Item(104, 15).Value = 2.8
$ws.Cells.Item(103, 16).Value = 0.75
$ws.Cells.Item(104, 16).Value = 0
$ws.Cells.Item(103, 17).Value = 1.95
$ws.Cells.Item(104, 17).Value = 1.75
$ws.Cells.Item(103, 18).Value = 1.75
$ws.Cells.Item(104, 18).Value = 2.05
$ws.Cells.Item(103, 19).Value = 3
$ws.Cells.Item(104, 19).Value = 2.25
$ws.Cells.Item(103, 20).Value = 1.975
$ws.Cells.Item(104, 20).Value = 1.875
$ws.Cells.Item(103, 21).Value = 1.825
$ws.Cells.Item(104, 21).Value = 1.925
$ws.Cells.Item(103, 22).Value = 3.2
$ws.Cells.Item(104, 22).Value = -1
$ws.Cells.Item(103, 23).Value = -1
$ws.Cells.Item(104, 23).Value = 2.1
$ws.Cells.Item(103, 25).Value = 0.95
$ws.Cells.Item(104, 25).Value = 0
$ws.Cells.Item(103, 26).Value = -1
$ws.Cells.Item(104, 26).Value = 0
$ws.Cells.Item(103, 27).Value = -1
$ws.Cells.Item(104, 27).Value = -0.5
$ws.Cells.Item(103, 28).Value = 0.825
$ws.Cells.Item(104, 28).Value = 0.4625
$ws.Cells.Item(107, 2).Value = 7874784
$ws.Cells.Item(108, 2).Value = 7874783
$ws.Cells.Item(107, 6).Value = "FK Atyrau"
$ws.Cells.Item(108, 6).Value = "FK Kyzylzhar"
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(108, 8).Value = 1
$ws.Cells.Item(107, 9).Value = "D"
$ws.Cells.Item(108, 9).Value = "H"
$ws.Cells.Item(107, 10).Value = 2.45
$ws.Cells.Item(108, 10).Value = 1.95
$ws.Cells.Item(107, 11).Value = 3
$ws.Cells.Item(108, 11).Value = 3.2
$ws.Cells.Item(107, 12).Value = 2.7
$ws.Cells.Item(108, 12).Value = 3.5
$ws.Cells.Item(107, 13).Value = 2.45
$ws.Cells.Item(108, 13).Value = 1.666
$ws.Cells.Item(107, 14).Value = 3
$ws.Cells.Item(108, 14).Value = 3.5
$ws.Cells.Item(107, 15).Value = 2.7
$ws.Cells.Item(108, 15).Value = 4.5
$ws.Cells.Item(107, 16).Value = 0
$ws.Cells.Item(108, 16).Value = -0.75
$ws.Cells.Item(107, 17).Value = 1.8
$ws.Cells.Item(108, 17).Value = 1.9
$ws.Cells.Item(107, 18).Value = 2
$ws.Cells.Item(108, 18).Value = 1.9
$ws.Cells.Item(107, 19).Value = 2.25
$ws.Cells.Item(108, 19).Value = 2.5
$ws.Cells.Item(107, 20).Value = 2
$ws.Cells.Item(108, 20).Value = 1.95
$ws.Cells.Item(107, 21).Value = 1.8
$ws.Cells.Item(108, 21).Value = 1.75
$ws.Cells.Item(107, 22).Value = -1
$ws.Cells.Item(108, 22).Value = 0.6659999999999999
$ws.Cells.Item(107, 23).Value = 2
$ws.Cells.Item(108, 23).Value = -1
$ws.Cells.Item(107, 25).Value = 0
$ws.Cells.Item(108, 25).Value = 0.45
$ws.Cells.Item(107, 26).Value = 0
$ws.Cells.Item(108, 26).Value = -0.5
$ws.Cells.Item(107, 27).Value = -1
$ws.Cells.Item(108, 27).Value = 0.95
$ws.Cells.Item(107, 28).Value = 0.8
$ws.Cells.Item(108, 28).Value = -1
$ws.Cells.Item(119, 2).Value = 7873759
$ws.Cells.Item(120, 2).Value = 7874795
$ws.Cells.Item(120, 5).Value = "Ordabasy"
$ws.Cells.Item(119, 6).Value = "FK Zhenys"
$ws.Cells.Item(120, 6).Value = "Tobol Kostanay"
$ws.Cells.Item(119, 7).Value = 3
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(119, 9).Value = "H"
$ws.Cells.Item(120, 9).Value = "D"
$ws.Cells.Item(119, 10).Value = 1.25
$ws.Cells.Item(120, 10).Value = 2.2
$ws.Cells.Item(119, 11).Value = 5.75
$ws.Cells.Item(120, 11).Value = 3.1
$ws.Cells.Item(119, 12).Value = 7
$ws.Cells.Item(120, 12).Value = 3
$ws.Cells.Item(119, 13).Value = 1.444
$ws.Cells.Item(120, 13).Value = 2.625
$ws.Cells.Item(119, 14).Value = 4.75
$ws.Cells.Item(120, 14).Value = 3
$ws.Cells.Item(119, 15).Value = 4.75
$ws.Cells.Item(120, 15).Value = 2.55
$ws.Cells.Item(119, 16).Value = -1.25
$ws.Cells.Item(120, 16).Value = 0
$ws.Cells.Item(119, 17).Value = 1.95
$ws.Cells.Item(120, 17).Value = 1.9
$ws.Cells.Item(119, 18).Value = 1.85
$ws.Cells.Item(120, 18).Value = 1.9
$ws.Cells.Item(119, 19).Value = 2.75
$ws.Cells.Item(120, 19).Value = 2
$ws.Cells.Item(119, 20).Value = 1.9
$ws.Cells.Item(120, 20).Value = 1.95
$ws.Cells.Item(119, 21).Value = 1.9
$ws.Cells.Item(120, 21).Value = 1.85
$ws.Cells.Item(119, 22).Value = 0.444
$ws.Cells.Item(120, 22).Value = -1
$ws.Cells.Item(119, 23).Value = -1
$ws.Cells.Item(120, 23).Value = 2
$ws.Cells.Item(119, 25).Value = 0.95
$ws.Cells.Item(120, 25).Value = 0
$ws.Cells.Item(119, 26).Value = -1
$ws.Cells.Item(120, 26).Value = 0
$ws.Cells.Item(119, 27).Value = 0.45
$ws.Cells.Item(120, 27).Value = -1
$ws.Cells.Item(119, 28).Value = -0.5
$ws.Cells.Item(120, 28).Value = 0.8500000000000001
$ws.Cells.Item(136, 2).Value = 7871218
$ws.Cells.Item(137, 2).Value = 7874807
$ws.Cells.Item(136, 6).Value = "FC Astana"
$ws.Cells.Item(137, 6).Value = "FK Kyzylzhar"
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(137, 7).Value = 2
$ws.Cells.Item(136, 8).Value = 2
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(136, 9).Value = "A"
$ws.Cells.Item(137, 9).Value = "H"
$ws.Cells.Item(136, 10).Value = 6.5
$ws.Cells.Item(137, 10).Value = 2.6
$ws.Cells.Item(136, 11).Value = 3.75
$ws.Cells.Item(137, 11).Value = 2.875
$ws.Cells.Item(136, 12).Value = 1.444
$ws.Cells.Item(137, 12).Value = 2.625
$ws.Cells.Item(136, 13).Value = 5
$ws.Cells.Item(137, 13).Value = 2.625
$ws.Cells.Item(136, 14).Value = 3.5
$ws.Cells.Item(137, 14).Value = 2.875
$ws.Cells.Item(136, 15).Value = 1.615
$ws.Cells.Item(137, 15).Value = 2.625
$ws.Cells.Item(136, 16).Value = 0.75
$ws.Cells.Item(137, 16).Value = 0
$ws.Cells.Item(136, 17).Value = 1.975
$ws.Cells.Item(137, 17).Value = 1.9
$ws.Cells.Item(136, 18).Value = 1.825
$ws.Cells.Item(137, 18).Value = 1.9
$ws.Cells.Item(136, 19).Value = 2.25
$ws.Cells.Item(137, 19).Value = 2
$ws.Cells.Item(136, 20).Value = 1.95
$ws.Cells.Item(137, 20).Value = 1.975
$ws.Cells.Item(136, 21).Value = 1.85
$ws.Cells.Item(137, 21).Value = 1.825
$ws.Cells.Item(136, 22).Value = -1
$ws.Cells.Item(137, 22).Value = 1.625
$ws.Cells.Item(136, 24).Value = 0.615
$ws.Cells.Item(137, 24).Value = -1
$ws.Cells.Item(136, 25).Value = -1
$ws.Cells.Item(137, 25).Value = 0.8999999999999999
$ws.Cells.Item(136, 26).Value = 0.825
$ws.Cells.Item(137, 26).Value = -1
$ws.Cells.Item(136, 27).Value = -0.5
$ws.Cells.Item(137, 27).Value = 0
$ws.Cells.Item(136, 28).Value = 0.425
$ws.Cells.Item(137, 28).Value = 0
